$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 20
$lastColOld = 11   # K

# --- Step 1: snapshot existing values (row 1..20, col A..K) before touching anything.
# Skip the known-blank "totalMean"/"totalStd" cells (old J/K, rows 3..20) entirely —
# merely reading an empty cell materializes it as an empty string, which we don't want.
$values = @{}
for ($r = 1; $r -le $lastRow; $r++) {
    for ($c = 1; $c -le $lastColOld; $c++) {
        if ($r -ge 3 -and ($c -eq 10 -or $c -eq 11)) { continue }
        $values["$r,$c"] = $ws.Cells.Item($r, $c).Value2
    }
}

# --- Step 2: shift every existing column one place to the right (K->L ... A->B), highest column first ---
for ($c = $lastColOld; $c -ge 1; $c--) {
    $destCol = $c + 1
    for ($r = 1; $r -le $lastRow; $r++) {
        if ($r -ge 3 -and ($c -eq 10 -or $c -eq 11)) { continue }
        $v = $values["$r,$c"]
        $destCell = $ws.Cells.Item($r, $destCol)
        if ($r -eq 1) {
            # header row: blank old A1 simply stays blank, nothing moves there
            if ($null -ne $v) { $destCell.Value = $v }
        } elseif ($c -eq 1) {
            # the segment-name column (old A) loses the bold/border style when it becomes column B
            if ($null -ne $v) { $destCell.Value = $v }
            $destCell.Style = "Normal"
        } else {
            if ($null -ne $v) { $destCell.Value = $v }
        }
    }
}

# Make sure the shifted-in totalMean/totalStd columns (K,L) stay genuinely blank for rows 3..20
for ($r = 3; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 11).ClearContents() | Out-Null
    $ws.Cells.Item($r, 12).ClearContents() | Out-Null
}

# L1 ("totalStd" header) is a brand-new cell (old sheet only went up to column K), so it
# doesn't inherit the bold/border/center header style the way B1..K1 do automatically — copy it in.
$ws.Range("K1").Copy() | Out-Null
$ws.Range("L1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Application.CutCopyMode = $false

# --- Step 3: new column A ---
# Header "segments" goes in B1 (keeps the bold/border/center style already sitting there); A1 stays blank.
$ws.Range("B1").Value = "segments"

# Segment index values 0..18 in A2:A20, styled like the header/name column (style index 1: bold, bordered, centered).
$ws.Range("C1").Copy() | Out-Null
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $r - 2
    $cell.PasteSpecial(-4122) | Out-Null   # xlPasteFormats — reuse the existing header/index style (index 1)
}
$ws.Application.CutCopyMode = $false
